# data cleanup continued in player_per_game_df
# Remove the "Leandro Barbosa" row from the player/award pivot table.
# Deleting the entire row (with a shift-up) automatically re-numbers the
# rows below it and drops the now-unused "Leandro Barbosa" shared string
# on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRow = 0
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 1).Value2 -eq "Leandro Barbosa") {
        $targetRow = $r
        break
    }
}

if ($targetRow -gt 0) {
    $ws.Rows.Item($targetRow).Delete()
}
